$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "72.573.89"
Set-TextValue "E2" "  +5.75%  "

Set-TextValue "D3" "2.645.49"
Set-TextValue "E3" "  +5.53%  "

Set-TextValue "E4" "  -0.06%  "

Set-TextValue "D5" "609.65"
Set-TextValue "E5" "  +3.17%  "

Set-TextValue "D6" "181.33"
Set-TextValue "E6" "  +4.05%  "

Set-TextValue "E7" "  -0.15%  "

Set-TextValue "D8" "0.528"
Set-TextValue "E8" "  +2.58%  "

Set-TextValue "D9" "0.177"
Set-TextValue "E9" "  +18.33%  "

Set-TextValue "D10" "2.644.13"
Set-TextValue "E10" "  +5.50%  "

Set-TextValue "E11" "  +1.00%  "

Set-TextValue "D12" "0.352"
Set-TextValue "E12" "  +5.12%  "

Set-TextValue "D13" "5.08"
Set-TextValue "E13" "  +1.42%  "

Set-TextValue "B14" "ShibaInu"
Set-TextValue "C14" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000197"
Set-TextValue "E14" "  +14.45%  "

Set-TextValue "B15" "WrappedliquidstakedEther2.0"
Set-TextValue "C15" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D15" "3.139.71"
Set-TextValue "E15" "  +6.45%  "

Set-TextValue "B16" "Avalanche"
Set-TextValue "C16" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D16" "26.86"
Set-TextValue "E16" "  +4.52%  "

Set-TextValue "B17" "WrappedBTC"
Set-TextValue "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "72.403.76"
Set-TextValue "E17" "  +5.52%  "

Set-TextValue "D18" "2.682.05"
Set-TextValue "E18" "  +7.26%  "

Set-TextValue "D19" "385.57"
Set-TextValue "E19" "  +6.33%  "

Set-TextValue "D20" "11.66"
Set-TextValue "E20" "  +6.97%  "

Set-TextValue "D21" "7.93"
Set-TextValue "E21" "  +5.31%  "

Set-TextValue "D22" "4.21"
Set-TextValue "E22" "  +4.95%  "

Set-TextValue "D23" "2.02"
Set-TextValue "E23" "  +22.39%  "

Set-TextValue "D24" "73.34"
Set-TextValue "E24" "  +4.60%  "

Set-TextValue "E25" "  +7.21%  "

Set-TextValue "D27" "10.02"
Set-TextValue "E27" "  +12.52%  "

Set-TextValue "D28" "2.782.64"

Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.48%  "

Set-TextValue "D30" "0.0₃0974"
Set-TextValue "E30" "  +11.17%  "

Set-TextValue "D31" "548.71"
Set-TextValue "E31" "  +7.36%  "

Set-TextValue "D32" "8.10"
Set-TextValue "E32" "  +5.05%  "

Set-TextValue "E33" "  +10.44%  "

Set-TextValue "E34" "  +4.35%  "

Set-TextValue "E35" "  -0.24%  "

Set-TextValue "D36" "166.23"
Set-TextValue "E36" "  +2.40%  "

Set-TextValue "B37" "EthereumClassic"
Set-TextValue "C37" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D37" "19.40"
Set-TextValue "E37" "  +4.76%  "

Set-TextValue "B38" "Kaspa"
Set-TextValue "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.114"
Set-TextValue "E38" "  -2.85%  "

Set-TextValue "D39" "1.42"
Set-TextValue "E39" "  +8.24%  "

Set-TextValue "D40" "19.12"
Set-TextValue "E40" "  +2.61%  "

Set-TextValue "D41" "1.87"
Set-TextValue "E41" "  +10.28%  "

Set-TextValue "D42" "5.13"
Set-TextValue "E42" "  +8.33%  "

Set-TextValue "D43" "2.63"
Set-TextValue "E43" "  +13.81%  "

Set-TextValue "E44" "  +0.13%  "

Set-TextValue "D45" "0.336"
Set-TextValue "E45" "  +6.10%  "

Set-TextValue "D46" "39.61"
Set-TextValue "E46" "  +1.83%  "

Set-TextValue "D47" "151.77"
Set-TextValue "E47" "  +1.22%  "

Set-TextValue "D48" "3.70"
Set-TextValue "E48" "  +4.42%  "

Set-TextValue "B49" "BabyDogeCoin"
Set-TextValue "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.0₆0272"
Set-TextValue "E49" "  +10.30%  "

Set-TextValue "B50" "ARBITRUM"
Set-TextValue "C50" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D50" "0.542"
Set-TextValue "E50" "  +6.04%  "

Set-TextValue "B51" "Optimism"
Set-TextValue "C51" "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextValue "D51" "1.71"
Set-TextValue "E51" "  +9.63%  "
